$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.940.88"
$ws.Range("E2").Value = "  +2.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.567.58"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.45"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.34"
$ws.Range("E6").Value = "  -3.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").Value = "  -5.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.577.06"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.014.94"
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.001.30"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.54"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.574.42"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.34"
$ws.Range("E20").Value = "  +4.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.28"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.03"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.87"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.687.28"
$ws.Range("E26").Value = "  +5.22%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0846"
$ws.Range("E29").Value = "  +4.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.42"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.97"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.15"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.55"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.70"
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.96"
$ws.Range("E36").Value = "  +3.00%  "
$ws.Range("E37").Value = "  +1.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  +23.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "298.44"
$ws.Range("E42").Value = "  +5.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.44"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0563"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0995"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.614"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.994"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.69"
$ws.Range("E48").Value = "  +8.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.94"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.029.73"
$ws.Range("E50").Value = "  +6.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0233"
$ws.Range("E51").Value = "  +0.06%  "
